# Re-order the (Fecha, Volumen, Precio min/max/promedio/$-kg) data across rows 2-13
# The underlying per-row records (D, M, N, O, P, S) get shuffled to new rows while
# all other columns (A, B, C, E-L, Q, R, T) stay put, since they are identical for
# every record in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for the columns that move, keyed by row number.
$rows = 2..13
$before = @{}
foreach ($r in $rows) {
    $rec = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        S = $ws.Cells.Item($r, 19).Value2
    }
    $before[$r] = $rec
}

# Mapping of destination row -> source row (where its new data comes from).
$mapping = @{
    2  = 12
    3  = 5
    4  = 6
    5  = 8
    6  = 4
    7  = 10
    8  = 7
    9  = 13
    10 = 9
    11 = 3
    12 = 11
    13 = 2
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $rec = $before[$srcRow]

    $ws.Cells.Item($destRow, 4).Value  = $rec.D
    $ws.Cells.Item($destRow, 13).Value = $rec.M
    $ws.Cells.Item($destRow, 14).Value = $rec.N
    $ws.Cells.Item($destRow, 15).Value = $rec.O
    $ws.Cells.Item($destRow, 16).Value = $rec.P
    $ws.Cells.Item($destRow, 19).Value = $rec.S
}
